$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("L7:L12").NumberFormat = "0.000%"
$ws.Range("L7").Formula = '=(E7-J7)/($E$12-$E$7)'
$ws.Range("L8:L12").Formula = '=(E8-J8)/($E$12-$E$7)'
